$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.800.84'
$ws.Range("E2").Value = '  -3.85%  '
$ws.Range("D3").Value = '1.817.62'
$ws.Range("E3").Value = '  -2.77%  '
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = "'277.73"
$ws.Range("E5").Value = '  -7.65%  '
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = "'0.5103"
$ws.Range("E7").Value = '  -4.49%  '
$ws.Range("D8").Value = "'0.3522"
$ws.Range("E8").Value = '  -6.16%  '
$ws.Range("D9").Value = "'45.38"
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").Value = "'0.06672"
$ws.Range("E10").Value = '  -6.99%  '
$ws.Range("D11").Value = "'19.98"
$ws.Range("E11").Value = '  -7.46%  '
$ws.Range("D12").Value = "'0.8337"
$ws.Range("E12").Value = '  -5.86%  '
$ws.Range("D13").Value = "'0.07914"
$ws.Range("E13").Value = '  -2.68%  '
$ws.Range("D14").Value = '1.823.52'
$ws.Range("E14").Value = '  -2.06%  '
$ws.Range("D15").Value = "'5.088"
$ws.Range("E15").Value = '  -3.33%  '
$ws.Range("D16").Value = "'87.81"
$ws.Range("E16").Value = '  -5.55%  '
$ws.Range("D17").Value = "'0.9989"
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("D18").Value = "'14.14"
$ws.Range("E18").Value = '  -3.96%  '
$ws.Range("D19").Value = "'0.000008025"
$ws.Range("E19").Value = '  -5.91%  '
$ws.Range("D20").Value = "'0.9995"
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").Value = '25.853.21'
$ws.Range("E21").Value = '  -3.52%  '
$ws.Range("E22").Value = '  -4.91%  '
$ws.Range("D23").Value = "'10.02"
$ws.Range("E23").Value = '  -6.22%  '
$ws.Range("D24").Value = "'6.081"
$ws.Range("E24").Value = '  -4.58%  '
$ws.Range("D25").Value = "'142.70"
$ws.Range("E25").Value = '  -2.98%  '
$ws.Range("D26").Value = "'2.190"
$ws.Range("E26").Value = '  -2.52%  '
$ws.Range("D27").Value = "'1.669"
$ws.Range("E27").Value = '  -3.56%  '
$ws.Range("D28").Value = "'17.12"
$ws.Range("E28").Value = '  -4.97%  '
$ws.Range("D29").Value = "'109.71"
$ws.Range("E29").Value = '  -4.07%  '
$ws.Range("D30").Value = "'4.336"
$ws.Range("E30").Value = '  -8.47%  '
$ws.Range("D31").Value = "'4.249"
$ws.Range("E31").Value = '  -7.11%  '
$ws.Range("D32").Value = "'0.08841"
$ws.Range("E32").Value = '  -3.19%  '
$ws.Range("D33").Value = "'0.04866"
$ws.Range("E33").Value = '  -2.46%  '
$ws.Range("D34").Value = "'0.7358"
$ws.Range("E34").Value = '  -8.07%  '
$ws.Range("D35").Value = "'1.138"
$ws.Range("E35").Value = '  -2.45%  '
$ws.Range("D36").Value = "'2.877"
$ws.Range("E36").Value = '  -3.97%  '
$ws.Range("D37").Value = "'3.156"
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("D38").Value = "'0.9989"
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").Value = "'0.5258"
$ws.Range("E39").Value = '  -11.20%  '
$ws.Range("D40").Value = "'2.326"
$ws.Range("E40").Value = '  -10.63%  '
$ws.Range("D41").Value = "'0.01849"
$ws.Range("E41").Value = '  -5.09%  '
$ws.Range("D42").Value = "'0.9617"
$ws.Range("E42").Value = '  -9.92%  '
$ws.Range("D43").Value = "'111.82"
$ws.Range("E43").Value = '  -2.96%  '
$ws.Range("D44").Value = "'6.195"
$ws.Range("E44").Value = '  -6.19%  '
$ws.Range("D45").Value = "'8.101"
$ws.Range("E45").Value = '  -8.55%  '
$ws.Range("D46").Value = "'0.9989"
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = "'0.4614"
$ws.Range("E47").Value = '  -8.93%  '
$ws.Range("D48").Value = "'0.1369"
$ws.Range("E48").Value = '  -8.32%  '
$ws.Range("E49").Value = '  -3.01%  '
$ws.Range("D50").Value = "'9.256"
$ws.Range("E50").Value = '  -6.60%  '
$ws.Range("D51").Value = "'1.503"
$ws.Range("E51").Value = '  -7.20%  '
